# Update "paises" (countries) workbook: refresh COVID stats snapshot and
# update the "last updated" timestamp. A handful of countries swapped rank
# (and therefore table row) because their updated totals crossed a
# neighboring country's total; for those pairs we write the country name
# together with its refreshed stats into the row it now occupies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp on row 1 -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 11 de Julio de 2020 a las 22:53"

# --- Helper data: row -> (Country, CasosTotales, NuevosCasos, CasosActivos,
#                           Recuperados, CasosCriticos, MuertesHoy, Muertes)
# Only rows whose stats changed are listed; country name is included so
# rows that were reordered (13/14, 131/132, 181/182) end up with the right
# name next to the right numbers.

$rows = @(
    @{ R=4;   A="Estados Unidos"; B=3343772; C=51986; D=1486312; E=1720148; F=0; G=641; H=137312 },
    @{ R=13;  A="Sudafrica";      B=264184;  C=13497; D=127715;  E=132498;  F=0; G=111; H=3971 },
    @{ R=14;  A="Iran";           B=255117;  C=2397;  D=217666;  E=24816;   F=0; G=188; H=12635 },
    @{ R=19;  A="Alemania";       B=199812;  C=224;   D=184500;  E=6178;    F=0; G=4;   H=9134 },
    @{ R=27;  A="Egipto";         B=81158;   C=923;   D=23876;   E=53513;   F=0; G=67;  H=3769 },
    @{ R=49;  A="Suiza";          B=32817;   C=127;   D=29500;   E=1349;    F=0; G=2;   H=1968 },
    @{ R=131; A="Ruanda";         B=1299;    C=47;    D=663;     E=632;     F=0; G=1;   H=4 },
    @{ R=132; A="Benin";          B=1285;    C=0;     D=333;     E=929;     F=0; G=0;   H=23 },
    @{ R=181; A="Monaco";         B=109;     C=1;     D=96;      E=9;       F=0; G=0;   H=4 },
    @{ R=182; A="Bahamas";        B=108;     C=1;     D=89;      E=8;       F=0; G=0;   H=11 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
}
